$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Bold the 4 existing collaborator names (rPr gains <w:b/>)
# ---------------------------------------------------------------------
$names = @("João Dinis", "Isaac Flores", "office 365 dev account", "Microsoft Open Source")
foreach ($name in $names) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($name, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1
    }
}

# ---------------------------------------------------------------------
# 2) Locate the 4 paragraphs that make up the 5th ("nam") collaborator
#    block, identified by their stable w14:paraId values.
# ---------------------------------------------------------------------
$paraRepo4   = $null   # 357CD8DD - "So repository co: 0" (4th person) - loses _GoBack bookmark
$paraName5   = $null   # 33ED76F8 - "Ten day du (nguoi thu nhat): " - becomes the 5th person's name
$paraGithub5 = $null   # 37F9D3B2 - "Link tai khoan Github : " - gets the github url appended
$paraRepo5   = $null   # 2FBB395B - "So repository co: " - gets the repo count appended

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $xml = $p.Range.WordOpenXML
    if ($xml.IndexOf('w14:paraId="357CD8DD"') -ge 0) { $paraRepo4 = $i }
    elseif ($xml.IndexOf('w14:paraId="33ED76F8"') -ge 0) { $paraName5 = $i }
    elseif ($xml.IndexOf('w14:paraId="37F9D3B2"') -ge 0) { $paraGithub5 = $i }
    elseif ($xml.IndexOf('w14:paraId="2FBB395B"') -ge 0) { $paraRepo5 = $i }
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 3) Paragraph 357CD8DD: drop the _GoBack bookmark (it moves to the
#    new paragraph below).
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item($paraRepo4)
$p1xml = '<w:p w14:paraId="357CD8DD" w14:textId="06EB066B" w:rsidR="0052705C" w:rsidRDefault="0052705C" w:rsidP="0052705C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Số repository có: </w:t></w:r><w:r w:rsidR="00DA2C77"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>0</w:t></w:r></w:p>'
$p1.Range.InsertXML($pkgHeader + $p1xml + $pkgFooter)

# ---------------------------------------------------------------------
# 4) Paragraph 33ED76F8: split "Tên đầy đủ (người thứ nhất): " into the
#    "Tên đầy đủ (người thứ năm): " phrasing, followed by the bolded
#    5th collaborator's name, and receives the relocated _GoBack
#    bookmark.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item($paraName5)
$p2xml = '<w:p w14:paraId="33ED76F8" w14:textId="77777777" w:rsidR="0052705C" w:rsidRDefault="0052705C" w:rsidP="0052705C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Tên đầy đủ (người thứ</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> năm </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>):</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:iCs/></w:rPr><w:t>Craig Presti</w:t></w:r><w:bookmarkStart w:id="10" w:name="_GoBack"/><w:bookmarkEnd w:id="10"/></w:p>'
$p2.Range.InsertXML($pkgHeader + $p2xml + $pkgFooter)

# ---------------------------------------------------------------------
# 5) Paragraph 37F9D3B2: append the github url as a new run.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item($paraGithub5)
$p3xml = '<w:p w14:paraId="37F9D3B2" w14:textId="77777777" w:rsidR="0052705C" w:rsidRDefault="0052705C" w:rsidP="0052705C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Link tài khoản Github : </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>https://github.com/craigomatic</w:t></w:r></w:p>'
$p3.Range.InsertXML($pkgHeader + $p3xml + $pkgFooter)

# ---------------------------------------------------------------------
# 6) Paragraph 2FBB395B: append the repository count as a new run.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item($paraRepo5)
$p4xml = '<w:p w14:paraId="2FBB395B" w14:textId="50E831EE" w:rsidR="0052705C" w:rsidRPr="00DA2C77" w:rsidRDefault="0052705C" w:rsidP="00DA2C77"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="40"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Số repository có: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>31</w:t></w:r></w:p>'
$p4.Range.InsertXML($pkgHeader + $p4xml + $pkgFooter)

Write-Output "Done"
